# Updated tasks for Jan
$wb = $excel.ActiveWorkbook

$wsJan = $wb.Worksheets.Item("January")
$wsFeb = $wb.Worksheets.Item("February")

# The team roster that used to live on the January sheet now also lives on
# the February sheet (B3:B6).
$names = @("Nitesh", "Gautami", "Pratiksha", "Pruthviraj")
for ($i = 0; $i -lt $names.Length; $i++) {
    $row = 3 + $i
    $wsFeb.Cells.Item($row, 2).Value = $names[$i]
}

# January now also tracks each person's task for the month in column C.
$tasks = @("Migration Testing", "Automation Testing", "API Testing", "Performance Testing")
for ($i = 0; $i -lt $tasks.Length; $i++) {
    $row = 3 + $i
    $wsJan.Cells.Item($row, 3).Value = $tasks[$i]
}

# Widen the new task column so the text fits.
$wsJan.Columns.Item(3).ColumnWidth = 18.7109375

# Add a new, blank sheet for April at the end of the workbook.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsApr = $wb.Worksheets.Add([Type]::Missing, $lastSheet)
$wsApr.Name = "April"

# Update the selections to cover the populated data ranges, and make
# February the active tab/sheet.
$wsJan.Range("B3:B6").Select() | Out-Null
$wsFeb.Activate() | Out-Null
$wsFeb.Range("B3:B6").Select() | Out-Null
